$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.857.25"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "1.802.56"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.15"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +4.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3702"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07360"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8683"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("D12").Value = "1.782.48"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.55"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.514"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07030"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008697"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.68"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "26.860.98"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.284"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("D24").Value = "2.056.90"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.906"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.34"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.139"
$ws.Range("E28").Value = "  -7.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.253"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.15"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08908"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7616"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.941"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.463"
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9998"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.100"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05249"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.225"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5306"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.370"
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("E44").Value = "  -3.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.510"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5020"
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.28"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.04"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9995"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.659"
$ws.Range("E50").Value = "  -2.30%  "
